# Regenerate merged AHB files
# - rename the "_old" / "_new" header columns to "_FV2404" / "_FV2410"
# - wrap the used range in a real Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = "_old" -> "_FV2404", L1:U1 = "_new" -> "_FV2410"; K1 "diff" stays) ---
$ws.Range("A1").Value2 = "Segmentname_FV2404"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2404"
$ws.Range("C1").Value2 = "Segment_FV2404"
$ws.Range("D1").Value2 = "Datenelement_FV2404"
$ws.Range("E1").Value2 = "Segment ID_FV2404"
$ws.Range("F1").Value2 = "Code_FV2404"
$ws.Range("G1").Value2 = "Qualifier_FV2404"
$ws.Range("H1").Value2 = "Beschreibung_FV2404"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value2 = "Bedingung_FV2404"

$ws.Range("L1").Value2 = "Segmentname_FV2410"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2410"
$ws.Range("N1").Value2 = "Segment_FV2410"
$ws.Range("O1").Value2 = "Datenelement_FV2410"
$ws.Range("P1").Value2 = "Segment ID_FV2410"
$ws.Range("Q1").Value2 = "Code_FV2410"
$ws.Range("R1").Value2 = "Qualifier_FV2410"
$ws.Range("S1").Value2 = "Beschreibung_FV2410"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value2 = "Bedingung_FV2410"

# --- 2. Turn the data range into an Excel Table ("Table1") ---
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U62"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row (split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
